$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = 14822.775
    3 = 16183.375
    4 = 7787.27
    5 = 10120.265
    6 = 12770.29
    7 = 10594.15
    8 = 15505.735
    9 = 12739.375
    10 = 13429.1
    11 = 23022.79
    12 = 12174.515
    13 = 10593.61
    14 = 10901.225
    15 = 7033.235000000001
    16 = 15737.085
    17 = 18498.915
    18 = 18643.575
    19 = 21642.04
    20 = 11046.62
    21 = 7348.895
    22 = 11146.28
    23 = 7108.69
    24 = 15275.505
    25 = 12856.035
    26 = 7911.139999999999
    27 = 11480.08
    28 = 11647.46
    29 = 12682.465
    30 = 24457.385
    31 = 14881
    32 = 6648.7
    33 = 13227.545
    34 = 9878.459999999999
    35 = 11863.175
    36 = 7092.425
    37 = 10446.95
    38 = 40740.27
    39 = 11898.695
    40 = 6557.974999999999
    41 = 15256.83
    42 = 11579.28
    43 = 10738.27
    44 = 13330.395
    45 = 10572.55
    46 = 16521.335
    47 = 12857.965
    48 = 4301.77
    49 = 14503.91
    50 = 6096.920000000001
    51 = 11006.865
    52 = 7935.965
    53 = 6112.22
    54 = 16366.12
    55 = 16703.545
    56 = 8268.924999999999
    57 = 9782.84
    58 = 16606.66
    59 = 9468.219999999999
    60 = 8930.870000000001
    61 = 13431.95
    62 = 12638.53
    63 = 9874.15
    64 = 10140.02
    65 = 15097.305
    66 = 9229.500000000002
    67 = 11720.15
    68 = 14485.04
    69 = 7233.25
    70 = 7915.825
    71 = 11342.025
    72 = 8662.519999999999
    73 = 11019.945
    74 = 12880.695
    75 = 11153.655
    76 = 9603.295
    77 = 17400.41
    78 = 10517.42
    79 = 15297.725
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
